$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Column C ("Förändrad") on rows 2..34: bump date serial 45643 -> 45644
#    (2024-12-17 -> 2024-12-18), preserving the existing date number format.
for ($r = 2; $r -le 34; $r++) {
    $ws.Cells.Item($r, 3).Value = 45644
}

# 2) Row 34 gains an explicit row height (15pt, custom height flag).
$ws.Rows.Item(34).RowHeight = 15

# 3) Append two new data rows (35 and 36) with the same shape as the
#    existing rows (A..E, G..Q values, R styled/wrapped empty cell).

# Row 35: A 60501-2024
$ws.Cells.Item(35, 1).Value = "A 60501-2024"
$ws.Cells.Item(35, 2).Value = 45643
$ws.Cells.Item(35, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(35, 3).Value = 45644
$ws.Cells.Item(35, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(35, 4).Value = "OKÄNT"
$ws.Cells.Item(35, 5).Value = "OKÄNT"
$ws.Cells.Item(35, 7).Value = 0.6
for ($col = 8; $col -le 17; $col++) {
    $ws.Cells.Item(35, $col).Value = 0
}
$ws.Cells.Item(35, 18).WrapText = $true
$ws.Rows.Item(35).RowHeight = 15

# Row 36: A 60500-2024
$ws.Cells.Item(36, 1).Value = "A 60500-2024"
$ws.Cells.Item(36, 2).Value = 45643
$ws.Cells.Item(36, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(36, 3).Value = 45644
$ws.Cells.Item(36, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(36, 4).Value = "OKÄNT"
$ws.Cells.Item(36, 5).Value = "OKÄNT"
$ws.Cells.Item(36, 7).Value = 0.8
for ($col = 8; $col -le 17; $col++) {
    $ws.Cells.Item(36, $col).Value = 0
}
$ws.Cells.Item(36, 18).WrapText = $true
